$wb = $excel.ActiveWorkbook

# Switch to the "Step 3" sheet and make it the active sheet/tab
$ws = $wb.Worksheets.Item("Step 3")
$ws.Activate()

# Enter the naive cost-calculation formula in K3 (mirrors K2/K3 pattern on other steps: H*G)
$ws.Range("K3").Formula = "=H3*G3"

# Leave selection on K5, as it ends up after entering the formula and moving down
$ws.Range("K5").Select()
